# Remove the trailing "Ver no Jupiter..." / copyright footer block that
# a Jekyll site rebuild dropped from the end of the bibliography section,
# along with the blank paragraph that separated it from the last
# reference entry.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph via Find so we don't depend
# on hard-coded paragraph numbers.
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if ($found) {
    # Figure out which paragraph index that Find hit corresponds to.
    $count = $d.Paragraphs.Count
    $targetIndex = -1
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -eq $findRange.Start) {
            $targetIndex = $i
        }
    }

    if ($targetIndex -gt 1) {
        # Remove the blank paragraph right before it, the "Ver no
        # Jupiter..." paragraph itself, and the following copyright
        # paragraph (target - 1 .. target + 1).
        $startPara = $d.Paragraphs.Item($targetIndex - 1)
        $endPara = $d.Paragraphs.Item($targetIndex + 1)

        $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
        $delRange.Delete()
    }
}
